$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1550
$ws.Range("I39").Value = 1100
$ws.Range("J39").Value = 2000
$ws.Range("K39").Value = 3300
$ws.Range("L39").Value = 6000
$ws.Range("M39").Value = -3004
$ws.Range("N39").Value = -6592

$ws.Range("H40").Value = 1216.3334
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1216.3334
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1216.3334
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -1566.3334

$ws.Range("H74").Value = 4403.091
$ws.Range("J74").Value = 4463.4
$ws.Range("L74").Value = 4463.4
$ws.Range("N74").Value = -6335.4

$ws.Range("H77").Value = 4403.091
$ws.Range("J77").Value = 4463.4
$ws.Range("L77").Value = 22317
$ws.Range("N77").Value = -31677

$ws.Range("H86").Value = 110537600
$ws.Range("I86").Value = 190927840
$ws.Range("J86").Value = 1025
$ws.Range("K86").Value = 190927840
$ws.Range("L86").Value = 1025
$ws.Range("M86").Value = -190926717
$ws.Range("N86").Value = -3271

$ws.Range("H89").Value = 110537600
$ws.Range("I89").Value = 190927840
$ws.Range("J89").Value = 1025
$ws.Range("K89").Value = 954639200
$ws.Range("L89").Value = 5125
$ws.Range("M89").Value = -954633584
$ws.Range("N89").Value = -16357

$ws.Range("H98").Value = 1053.9
$ws.Range("I98").Value = 837.6667
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 837.6667
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = 660.3333
$ws.Range("N98").Value = -5996

$ws.Range("H112").Value = 6475.8335
$ws.Range("I112").Value = 1550
$ws.Range("J112").Value = 6690
$ws.Range("K112").Value = 4650
$ws.Range("L112").Value = 20070
$ws.Range("M112").Value = -3542
$ws.Range("N112").Value = -22286

$ws.Range("H113").Value = 2501.8
$ws.Range("I113").Value = 2170
$ws.Range("J113").Value = 2999.5
$ws.Range("K113").Value = 2170
$ws.Range("L113").Value = 2999.5
$ws.Range("M113").Value = 1084
$ws.Range("N113").Value = -9507.5

$ws.Range("H122").Value = 1053.9
$ws.Range("I122").Value = 837.6667
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 2513.0001
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -63.0001000000002
$ws.Range("N122").Value = -13900

$ws.Range("H135").Value = 3666.36
$ws.Range("I135").Value = 2264.5
$ws.Range("J135").Value = 7271.143
$ws.Range("K135").Value = 20380.5
$ws.Range("L135").Value = 65440.287
$ws.Range("M135").Value = -17845.5
$ws.Range("N135").Value = -70510.287

$ws.Range("H138").Value = 1989.8541
$ws.Range("I138").Value = 2829.5334
$ws.Range("J138").Value = 1834.358
$ws.Range("K138").Value = 8488.600199999999
$ws.Range("L138").Value = 5503.074
$ws.Range("M138").Value = -3348.600199999999
$ws.Range("N138").Value = -15783.074


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3103.3044
$ws.Range("I61").Value = 2482.7693
$ws.Range("J61").Value = 3910
$ws.Range("K61").Value = 2482.7693
$ws.Range("L61").Value = 3910
$ws.Range("M61").Value = -2270.7693
$ws.Range("N61").Value = -4334

$ws.Range("H136").Value = 3103.3044
$ws.Range("I136").Value = 2482.7693
$ws.Range("J136").Value = 3910
$ws.Range("K136").Value = 7448.3079
$ws.Range("L136").Value = 11730
$ws.Range("M136").Value = -4898.3079
$ws.Range("N136").Value = -16830


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 47621720
$ws.Range("I86").Value = 58826300
$ws.Range("J86").Value = 2243.5
$ws.Range("K86").Value = 58826300
$ws.Range("L86").Value = 2243.5
$ws.Range("M86").Value = -58825177
$ws.Range("N86").Value = -4489.5

$ws.Range("H89").Value = 47621720
$ws.Range("I89").Value = 58826300
$ws.Range("J89").Value = 2243.5
$ws.Range("K89").Value = 294131500
$ws.Range("L89").Value = 11217.5
$ws.Range("M89").Value = -294125884
$ws.Range("N89").Value = -22449.5

$ws.Range("H94").Value = 877.38464
$ws.Range("I94").Value = 607.3333
$ws.Range("J94").Value = 1485
$ws.Range("K94").Value = 607.3333
$ws.Range("L94").Value = 1485
$ws.Range("M94").Value = -156.3333
$ws.Range("N94").Value = -2387


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2525000
$ws.Range("I6").Value = 2525000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 2525000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -2524887
$ws.Range("N6").ClearContents()


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 735.8095
$ws.Range("I5").Value = 649.05554
$ws.Range("K5").Value = 1947.16662
$ws.Range("M5").Value = -1835.16662

$ws.Range("H48").Value = 8400.4
$ws.Range("I48").Value = 1000
$ws.Range("J48").Value = 10250.5
$ws.Range("K48").Value = 3000
$ws.Range("L48").Value = 30751.5
$ws.Range("M48").Value = -2750
$ws.Range("N48").Value = -31251.5

$ws.Range("H51").Value = 1080.7693
$ws.Range("I51").Value = 580
$ws.Range("J51").Value = 1200
$ws.Range("K51").Value = 1740
$ws.Range("L51").Value = 3600
$ws.Range("M51").Value = -1280
$ws.Range("N51").Value = -4520

$ws.Range("H135").Value = 735.8095
$ws.Range("I135").Value = 649.05554
$ws.Range("K135").Value = 5841.49986
$ws.Range("M135").Value = -3306.49986

$ws.Range("H137").Value = 12047.154
$ws.Range("J137").Value = 4061.4
$ws.Range("L137").Value = 12184.2
$ws.Range("N137").Value = -22384.2

$ws.Range("H140").Value = 1418.25
$ws.Range("I140").Value = 1122.0952
$ws.Range("J140").Value = 1983.6364
$ws.Range("K140").Value = 3366.2856
$ws.Range("L140").Value = 5950.9092
$ws.Range("M140").Value = 1813.7144
$ws.Range("N140").Value = -16310.9092


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 943294.1
$ws.Range("J80").Value = 69161.164
$ws.Range("L80").Value = 69161.164
$ws.Range("N80").Value = -71157.164

$ws.Range("H83").Value = 943294.1
$ws.Range("J83").Value = 69161.164
$ws.Range("L83").Value = 345805.82
$ws.Range("N83").Value = -355789.82

$ws.Range("H122").Value = 4456.579
$ws.Range("I122").Value = 1590.909
$ws.Range("K122").Value = 4772.727000000001
$ws.Range("M122").Value = -2322.727000000001

$ws.Range("H132").Value = 3506
$ws.Range("I132").Value = 3119.111
$ws.Range("K132").Value = 9357.332999999999
$ws.Range("M132").Value = -6827.332999999999


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 9802.083000000001
$ws.Range("I22").Value = 1524.2858
$ws.Range("J22").Value = 21391
$ws.Range("K22").Value = 1524.2858
$ws.Range("L22").Value = 21391
$ws.Range("M22").Value = -1229.2858
$ws.Range("N22").Value = -21981

$ws.Range("H27").Value = 9802.083000000001
$ws.Range("I27").Value = 1524.2858
$ws.Range("J27").Value = 21391
$ws.Range("K27").Value = 1524.2858
$ws.Range("L27").Value = 21391
$ws.Range("M27").Value = -1417.2858
$ws.Range("N27").Value = -21605

$ws.Range("H40").Value = 2363.25
$ws.Range("I40").Value = 1982
$ws.Range("J40").Value = 2897
$ws.Range("K40").Value = 1982
$ws.Range("L40").Value = 2897
$ws.Range("M40").Value = -1846
$ws.Range("N40").Value = -3169

$ws.Range("H93").Value = 9392.134
$ws.Range("I93").Value = 12188.3
$ws.Range("J93").Value = 3799.8
$ws.Range("K93").Value = 12188.3
$ws.Range("L93").Value = 3799.8
$ws.Range("M93").Value = -10940.3
$ws.Range("N93").Value = -6295.8

$ws.Range("H132").Value = 4888.3076
$ws.Range("I132").Value = 4222
$ws.Range("J132").Value = 5665.6665
$ws.Range("K132").Value = 12666
$ws.Range("L132").Value = 16996.9995
$ws.Range("M132").Value = -10136
$ws.Range("N132").Value = -22056.9995


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 42400
$ws.Range("I70").Value = 13095
$ws.Range("J70").Value = 57052.5
$ws.Range("K70").Value = 13095
$ws.Range("L70").Value = 57052.5
$ws.Range("M70").Value = -12780
$ws.Range("N70").Value = -57682.5

$ws.Range("H73").Value = 42400
$ws.Range("I73").Value = 13095
$ws.Range("J73").Value = 57052.5
$ws.Range("K73").Value = 13095
$ws.Range("L73").Value = 57052.5
$ws.Range("M73").Value = -12003
$ws.Range("N73").Value = -59236.5

$ws.Range("H132").Value = 6175779
$ws.Range("I132").Value = 5212.143
$ws.Range("J132").Value = 8335477.5
$ws.Range("K132").Value = 15636.429
$ws.Range("L132").Value = 25006432.5
$ws.Range("M132").Value = -13106.429
$ws.Range("N132").Value = -25011492.5

$ws.Range("H136").Value = 2246.5173
$ws.Range("I136").Value = 1963
$ws.Range("K136").Value = 5889
$ws.Range("M136").Value = -3339

